# Update the "time_taken" column (F) on the "data" sheet with refreshed timestamps
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$data.Range("F2").Value = "2021-10-05 14:21:35.580499"
$data.Range("F3").Value = "2021-10-05 14:21:35.580507"
$data.Range("F4").Value = "2021-10-05 14:21:35.580511"
$data.Range("F5").Value = "2021-10-05 14:21:35.580513"
$data.Range("F6").Value = "2021-10-05 14:21:35.580517"
$data.Range("F7").Value = "2021-10-05 14:21:35.580519"
$data.Range("F8").Value = "2021-10-05 14:21:35.580522"
$data.Range("F9").Value = "2021-10-05 14:21:35.580525"
$data.Range("F10").Value = "2021-10-05 14:21:35.580527"
$data.Range("F11").Value = "2021-10-05 14:21:35.580530"
$data.Range("F12").Value = "2021-10-05 14:21:35.580533"

# Add a new "metadata" sheet after "data"
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "metadata"
$newSheet.Move($null, $wb.Worksheets.Item("data"))

$meta = $wb.Worksheets.Item("metadata")

# Header row (bold, matches the "data" sheet header style: thin border, center/top align)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$headerRng = $meta.Range("B1:G1")
$headerRng.Font.Bold = $true
$headerRng.Borders.LineStyle = 1
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160

# Data row
$meta.Range("A2").Value = 0
$a2 = $meta.Range("A2")
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$meta.Range("B2").Value = "Mitochondrial liver disease"
$meta.Range("C2").Value = 532
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.4"
$meta.Range("E2").Value = "2020-11-16T15:37:26.461851Z"
$meta.Range("F2").Value = "2021-10-05 14:21:35.577290"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/532/?format=json"
